$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.700.51'
$ws.Range("E2").Value = '  -4.20%  '
$ws.Range("D3").Value = '3.323.73'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.58'
$ws.Range("E5").Value = '  -3.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.42'
$ws.Range("E6").Value = '  -5.34%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("E9").Value = '  -4.14%  '
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.404'
$ws.Range("E11").Value = '  -4.70%  '
$ws.Range("D12").Value = '3.901.20'
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.18'
$ws.Range("E14").Value = '  -5.46%  '
$ws.Range("D15").Value = '66.746.60'
$ws.Range("E15").Value = '  -4.20%  '
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("D17").Value = '3.337.22'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.70'
$ws.Range("E19").Value = '  -2.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '432.70'
$ws.Range("E20").Value = '  -4.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.64'
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.58'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.08'
$ws.Range("E27").Value = '  -5.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.82'
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.33'
$ws.Range("E31").Value = '  -5.41%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.72'
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.35'
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("E38").Value = '  -4.31%  '
$ws.Range("D39").Value = '2.832.92'
$ws.Range("E39").Value = '  +3.14%  '
$ws.Range("E40").Value = '  -3.98%  '
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.20'
$ws.Range("E42").Value = '  -4.92%  '
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.52'
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("E46").Value = '  -7.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '324.27'
$ws.Range("E47").Value = '  -5.72%  '
$ws.Range("E48").Value = '  -4.80%  '
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("E50").Value = '  -2.50%  '
$ws.Range("E51").Value = '  -1.53%  '
